$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the "Roux Institute" distance column (O) values for data rows 2-17,
# leaving the existing cell style/formatting intact.
$ws.Range("O2:O17").ClearContents()
